# Auto-generated from the OOXML unified diff.
# All source cells in columns B:E are stored as text (inlineStr) in the
# original workbook, so every write below forces the destination cell to
# Text format first -- this preserves values like "1.00" / "20.10" / "3.30"
# that Excel would otherwise silently renormalise into numbers (dropping
# trailing zeros) when assigned through .Value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.941.74"
$ws.Range("E2").Value = "  +4.47%  "
$ws.Range("D3").Value = "2.679.39"
$ws.Range("E3").Value = "  +7.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.76"
$ws.Range("E5").Value = "  +9.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "326.06"
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.529"
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +3.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.95"
$ws.Range("E10").Value = "  +6.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.10"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("E12").Value = "  +3.37%  "
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("E14").Value = "  +5.10%  "
$ws.Range("D15").Value = "3.098.83"
$ws.Range("E15").Value = "  +8.06%  "
$ws.Range("D16").Value = "2.678.91"
$ws.Range("E16").Value = "  +6.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.875"
$ws.Range("E17").Value = "  +6.42%  "
$ws.Range("D18").Value = "49.900.92"
$ws.Range("E18").Value = "  +4.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("E19").Value = "  +4.07%  "
$ws.Range("E20").Value = "  +4.41%  "
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.19"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.89"
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.87"
$ws.Range("E26").Value = "  +4.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.14"
$ws.Range("E28").Value = "  +6.04%  "
$ws.Range("E29").Value = "  +3.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.27"
$ws.Range("E30").Value = "  +5.52%  "
$ws.Range("E31").Value = "  +3.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.31"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.49"
$ws.Range("E33").Value = "  +4.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.53"
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("E35").Value = "  +5.80%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.05"
$ws.Range("E37").Value = "  +12.45%  "
$ws.Range("E38").Value = "  +7.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.15"
$ws.Range("E39").Value = "  +10.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.84"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.84"
$ws.Range("E41").Value = "  +5.71%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.113"
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("E44").Value = "  +6.37%  "
$ws.Range("D45").Value = "2.117.20"
$ws.Range("E45").Value = "  +6.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.30"
$ws.Range("E46").Value = "  +6.06%  "
$ws.Range("E47").Value = "  +15.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.07"
$ws.Range("E48").Value = "  +9.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.06"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.36"
$ws.Range("E50").Value = "  +5.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.65"
$ws.Range("E51").Value = "  +7.13%  "
